$d = $word.ActiveDocument

# The paragraph currently reads "Version 1." and needs to become "Version 2."
# but with the runs restructured to match the target OOXML:
#   "Versi" | "on" | spellEnd | " 2" | bookmarkStart | bookmarkEnd | "."
#
# Apply edits back-to-front so earlier character offsets stay valid.

# 1) " 1." (positions 8-10, i.e. the "1" and the ".") becomes " 2"
#    (drop the trailing period here; it gets re-inserted after the bookmark below).
$tail = $d.Range(8, 10)
$tail.Text = "2"

# 2) Re-insert the "." as its own run, placed after the (hidden) _GoBack
#    bookmark that wraps the end of the paragraph's text.
$paraEnd = $d.Paragraphs(1).Range.End - 1
$period = $d.Range($paraEnd, $paraEnd)
$period.InsertAfter(".")

# 3) Split "Version" into two runs, "Versi" and "on", with no separator
#    between them (matching the target markup). A Font property is
#    toggled on/off to force the run boundary at offset 5.
$head = $d.Range(0, 5)
$head.Font.Bold = 1
$head.Font.Bold = 0

Write-Output $d.Content.Text
